$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet after the "url" sheet and name it "addSingleUser"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "addSingleUser"

# Populate the new sheet with data
$ws2.Range("A1").Value = "userRole"
$ws2.Range("B1").Value = "Admin"
$ws2.Range("C1").Value = "Admin"

$ws2.Range("A2").Value = "status"
$ws2.Range("B2").Value = "Enabled"
$ws2.Range("C2").Value = "Enabled"

$ws2.Range("A3").Value = "newUsername"
$ws2.Range("B3").Value = "tin01ak"
$ws2.Range("C3").Value = "tin02ak"

$ws2.Range("A4").Value = "newPassword"
$ws2.Range("B4").Value = "re5tr1ct$"
$ws2.Range("C4").Value = "re5tr1ct$"

# Select cell B7 on sheet1 (no longer the tab-selected sheet)
$ws1.Range("B7").Select()

# Select cell H21 and make addSingleUser the active sheet/tab
$ws2.Activate()
$ws2.Range("H21").Select()
